$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# --- 1. Duplicate the existing "2022-Q2" sheet so its original data/format
#        survives unchanged under the same name, then repurpose the
#        original tab as the new "2022-Q3" quarter. ---
$wsQ2.Copy($null, $wsQ2)
$wsQ2.Name = "2022-Q3"
$wsQ2Copy = $wb.Worksheets.Item(3)
$wsQ2Copy.Name = "2022-Q2"

# --- 2. Update the "总计" (summary) sheet: row 2 becomes the new quarter,
#        and a new row 3 repeats the old 2022-Q2 totals. ---
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q2"
$wsTotal.Cells.Item(3, 3).Value = 2
$wsTotal.Cells.Item(3, 4).Value = 0.03
$wsTotal.Cells.Item(2, 1).Copy()
$wsTotal.Cells.Item(3, 1).PasteSpecial(-4122)

# --- 3. Replace the "2022-Q3" sheet's holdings with the new quarter data. ---
# Re-stripe the header row + index column with the sheet's "header" style
# (matches the style already used by the 总计 sheet's own header/index cells).
$wsTotal.Cells.Item(2, 1).Copy()
$wsQ2.Range("B1:H1").PasteSpecial(-4122)
$wsQ2.Cells.Item(2, 1).PasteSpecial(-4122)
$wsQ2.Cells.Item(3, 1).PasteSpecial(-4122)

# Values that look numeric must stay text, like the source data -- force
# Text format before writing so Excel doesn't silently coerce them.
$wsQ2.Range("B2:B3").NumberFormat = "@"
$wsQ2.Range("D2:G3").NumberFormat = "@"

$wsQ2.Cells.Item(1, 2).Value = "基金代码"
$wsQ2.Cells.Item(1, 3).Value = "基金名称"
$wsQ2.Cells.Item(1, 4).Value = "基金规模"
$wsQ2.Cells.Item(1, 5).Value = "股票总仓位"
$wsQ2.Cells.Item(1, 6).Value = "仓位占比"
$wsQ2.Cells.Item(1, 7).Value = "持有市值(亿元)"
$wsQ2.Cells.Item(1, 8).Value = "仓位排名"

$wsQ2.Cells.Item(2, 1).Value = 0
$wsQ2.Cells.Item(2, 2).Value = "166109"
$wsQ2.Cells.Item(2, 3).Value = "信澳量化先锋混合（LOF）A"
$wsQ2.Cells.Item(2, 4).Value = "0.79"
$wsQ2.Cells.Item(2, 5).Value = "88.99"
$wsQ2.Cells.Item(2, 6).Value = "3.19"
$wsQ2.Cells.Item(2, 7).Value = "0.0252"
$wsQ2.Cells.Item(2, 8).Value = 6

$wsQ2.Cells.Item(3, 1).Value = 1
$wsQ2.Cells.Item(3, 2).Value = "166110"
$wsQ2.Cells.Item(3, 3).Value = "信澳量化先锋混合（LOF）C"
$wsQ2.Cells.Item(3, 4).Value = "0.11"
$wsQ2.Cells.Item(3, 5).Value = "88.99"
$wsQ2.Cells.Item(3, 6).Value = "3.19"
$wsQ2.Cells.Item(3, 7).Value = "0.0035"
$wsQ2.Cells.Item(3, 8).Value = 6

# The Text-format nudge above leaves those cells tagged with a non-default
# style; re-stamp them with the sheet's untouched default (A1) format so
# only the cell TYPE (text) sticks, not a stray numFmt/style.
$wsQ2.Cells.Item(1, 1).Copy()
$wsQ2.Range("B2:B3").PasteSpecial(-4122)
$wsQ2.Range("D2:G3").PasteSpecial(-4122)
